$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 4333.3335
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 62
$ws.Range("H62").Value = 5256.952
$ws.Range("I62").Value = 3769.5386
$ws.Range("K62").Value = 3769.5386
$ws.Range("M62").Value = -3145.5386

# Row 64
$ws.Range("H64").Value = 3284.8
$ws.Range("J64").Value = 3467.3
$ws.Range("L64").Value = 3467.3
$ws.Range("N64").Value = -3963.3

# Row 65
$ws.Range("H65").Value = 5256.952
$ws.Range("I65").Value = 3769.5386
$ws.Range("K65").Value = 18847.693
$ws.Range("M65").Value = -15727.693

# Row 67
$ws.Range("H67").Value = 3284.8
$ws.Range("J67").Value = 3467.3
$ws.Range("L67").Value = 3467.3
$ws.Range("N67").Value = -5183.3

# Row 98
$ws.Range("H98").Value = 755.63635
$ws.Range("I98").Value = 488
$ws.Range("J98").Value = 1142.2222
$ws.Range("K98").Value = 488
$ws.Range("L98").Value = 1142.2222
$ws.Range("M98").Value = 1010
$ws.Range("N98").Value = -4138.2222

# Row 122
$ws.Range("H122").Value = 755.63635
$ws.Range("I122").Value = 488
$ws.Range("J122").Value = 1142.2222
$ws.Range("K122").Value = 1464
$ws.Range("L122").Value = 3426.6666
$ws.Range("M122").Value = 986
$ws.Range("N122").Value = -8326.6666

# Row 129
$ws.Range("H129").Value = 173450.73
$ws.Range("J129").Value = 179636.53
$ws.Range("L129").Value = 538909.59
$ws.Range("N129").Value = -548909.59

# Row 137
$ws.Range("H137").Value = 1895.5454
$ws.Range("I137").Value = 2240.2
$ws.Range("J137").Value = 1608.3334
$ws.Range("K137").Value = 6720.599999999999
$ws.Range("L137").Value = 4825.0002
$ws.Range("M137").Value = -4170.599999999999
$ws.Range("N137").Value = -9925.0002

# Row 141
$ws.Range("H141").Value = 3397.7856
$ws.Range("I141").Value = 2960.818
$ws.Range("K141").Value = 8882.454000000002
$ws.Range("M141").Value = -3702.454000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 1954964.5
$ws.Range("I63").Value = 1979.0834
$ws.Range("J63").Value = 7813920.5
$ws.Range("K63").Value = 1979.0834
$ws.Range("L63").Value = 7813920.5
$ws.Range("M63").Value = -1293.0834
$ws.Range("N63").Value = -7815292.5

# Row 66
$ws.Range("H66").Value = 1954964.5
$ws.Range("I66").Value = 1979.0834
$ws.Range("J66").Value = 7813920.5
$ws.Range("K66").Value = 9895.416999999999
$ws.Range("L66").Value = 39069602.5
$ws.Range("M66").Value = -6463.416999999999
$ws.Range("N66").Value = -39076466.5

# Row 76
$ws.Range("H76").Value = 19999.5
$ws.Range("J76").Value = 19999.5
$ws.Range("L76").Value = 19999.5
$ws.Range("N76").Value = -20675.5

# Row 79
$ws.Range("H79").Value = 19999.5
$ws.Range("J79").Value = 19999.5
$ws.Range("L79").Value = 19999.5
$ws.Range("N79").Value = -22339.5

# Row 92
$ws.Range("H92").Value = 17592
$ws.Range("J92").Value = 17592
$ws.Range("L92").Value = 17592
$ws.Range("N92").Value = -22584

# Row 114
$ws.Range("H114").Value = 29362
$ws.Range("J114").Value = 29362
$ws.Range("L114").Value = 29362
$ws.Range("N114").Value = -38040

# Row 132
$ws.Range("H132").Value = 25390.182
$ws.Range("I132").Value = 2367
$ws.Range("K132").Value = 7101
$ws.Range("M132").Value = -4571

# Row 139
$ws.Range("H139").Value = 36136.637
$ws.Range("J139").Value = 36136.637
$ws.Range("L139").Value = 36136.637
$ws.Range("N139").Value = -46416.637

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 3707
$ws.Range("I94").Value = 1049.6666
$ws.Range("J94").Value = 5700
$ws.Range("K94").Value = 1049.6666
$ws.Range("L94").Value = 5700
$ws.Range("M94").Value = -598.6666
$ws.Range("N94").Value = -6602

# Row 115
$ws.Range("H115").Value = 27500
$ws.Range("J115").Value = 27500
$ws.Range("L115").Value = 27500
$ws.Range("N115").Value = -29850

# Row 132
$ws.Range("H132").Value = 2595.7812
$ws.Range("I132").Value = 1828.125
$ws.Range("J132").Value = 4898.75
$ws.Range("K132").Value = 5484.375
$ws.Range("L132").Value = 14696.25
$ws.Range("M132").Value = -2954.375
$ws.Range("N132").Value = -19756.25

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1515.6316
$ws.Range("I5").Value = 1345.4
$ws.Range("J5").Value = 1704.7778
$ws.Range("K5").Value = 4036.2
$ws.Range("L5").Value = 5114.3334
$ws.Range("M5").Value = -3924.2
$ws.Range("N5").Value = -5338.3334

# Row 68
$ws.Range("H68").Value = 10065.75
$ws.Range("J68").Value = 11778.9
$ws.Range("L68").Value = 35336.7
$ws.Range("N68").Value = -36958.7

# Row 71
$ws.Range("H71").Value = 10065.75
$ws.Range("J71").Value = 11778.9
$ws.Range("L71").Value = 106010.1
$ws.Range("N71").Value = -114122.1

# Row 86
$ws.Range("H86").Value = 35715040
$ws.Range("I86").Value = 750.2222
$ws.Range("J86").Value = 100000760
$ws.Range("K86").Value = 2250.6666
$ws.Range("L86").Value = 300002280
$ws.Range("M86").Value = -1064.6666
$ws.Range("N86").Value = -300004652

# Row 89
$ws.Range("H89").Value = 35715040
$ws.Range("I89").Value = 750.2222
$ws.Range("J89").Value = 100000760
$ws.Range("K89").Value = 6751.999800000001
$ws.Range("L89").Value = 900006840
$ws.Range("M89").Value = -823.9998000000005
$ws.Range("N89").Value = -900018696

# Row 131
$ws.Range("H131").Value = 749.7041
$ws.Range("J131").Value = 750.73193
$ws.Range("L131").Value = 2252.19579
$ws.Range("N131").Value = -12332.19579

# Row 135
$ws.Range("H135").Value = 1515.6316
$ws.Range("I135").Value = 1345.4
$ws.Range("J135").Value = 1704.7778
$ws.Range("K135").Value = 12108.6
$ws.Range("L135").Value = 15343.0002
$ws.Range("M135").Value = -9573.6
$ws.Range("N135").Value = -20413.0002

# Row 137
$ws.Range("H137").Value = 18524690
$ws.Range("I137").Value = 1265
$ws.Range("J137").Value = 20840118
$ws.Range("K137").Value = 3795
$ws.Range("L137").Value = 62520354
$ws.Range("M137").Value = 1305
$ws.Range("N137").Value = -62530554

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 26318918
$ws.Range("I102").Value = 35717440
$ws.Range("K102").Value = 35717440
$ws.Range("M102").Value = -35715818

# Row 116
$ws.Range("H116").Value = 41250
$ws.Range("J116").Value = 41250
$ws.Range("L116").Value = 41250
$ws.Range("N116").Value = -50428

# Row 132
$ws.Range("H132").Value = 44725.918
$ws.Range("I132").Value = 3001.5
$ws.Range("K132").Value = 9004.5
$ws.Range("M132").Value = -6474.5

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 805575.3
$ws.Range("I132").Value = 1507219
$ws.Range("K132").Value = 4521657
$ws.Range("M132").Value = -4519127

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3499.6667
$ws.Range("J62").Value = 3624.75
$ws.Range("L62").Value = 3624.75
$ws.Range("N62").Value = -4872.75

# Row 65
$ws.Range("H65").Value = 3499.6667
$ws.Range("J65").Value = 3624.75
$ws.Range("L65").Value = 18123.75
$ws.Range("N65").Value = -24363.75

# Row 126
$ws.Range("H126").Value = 949.24243
$ws.Range("I126").Value = 985.0714
$ws.Range("J126").Value = 748.6
$ws.Range("K126").Value = 2955.2142
$ws.Range("L126").Value = 2245.8
$ws.Range("M126").Value = -485.2142000000003
$ws.Range("N126").Value = -7185.8

# Row 132
$ws.Range("H132").Value = 1415.5312
$ws.Range("I132").Value = 899.9524
$ws.Range("J132").Value = 2399.818
$ws.Range("K132").Value = 2699.8572
$ws.Range("L132").Value = 7199.454000000001
$ws.Range("M132").Value = -169.8571999999999
$ws.Range("N132").Value = -12259.454
